$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190185546875
$ws.Range("B1").Value = 2.095541000366211
$ws.Range("C1").Value = 5.727797031402588
$ws.Range("D1").Value = 0.9269390106201172
$ws.Range("E1").Value = 1.104873418807983
